# Generate Report for Handoff
#
# The "b.md" row (row 3) on every sheet moves from the stale
# "Handed back: in sync with en-US" / old handoff-file state to a fresh
# "Ready for handoff" status, stamped with the new handoff datetime and
# pointing at the freshly generated handoff package
# (b.63290e5768f688058c7b37413b0a5c26c308f864.<locale>.xlf).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# Row 3 = b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-27-20 20:27:56"

# ---------------------------------------------------------------------
# zh-cn sheet detail row for b.md (row 3):
#   Status (C)              -> Ready for handoff
#   Latest Handoff File (D) -> new handoff package file name
#   Latest Handoff Datetime (E) -> new handoff timestamp
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-20 20:27:53"

# ---------------------------------------------------------------------
# de-de sheet detail row for b.md (row 3):
#   Status (C)              -> Ready for handoff
#   Latest Handoff File (D) -> new handoff package file name
#   Latest Handoff Datetime (E) -> new handoff timestamp
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-20 20:27:56"
